# edit.ps1 - Apply LOT2007 disciplina sheet restructuring
# Rebuilds rows 10-24 content per the target diff, then removes the
# now-redundant row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Long reused text blocks ----
$txtAdriane = "427823 - Adriane Maria Ferreira Milagres"
$txtShortSyllabusPt = "1.Acids and bases/Buffer solutions; 2. Amino acids; 3. Proteins: The primary level of protein structure; 4. Proteins: The three-dimensional structure; 5. Protein Function; 6. Enzymes: Biological catalysts; 7. The kinetics of enzymatic catalysis; 8. Carbohydrates; 9. Lipids; 10. Membranes and cellular transport; 11. Nucleic acids."
$txtSyllabusEn = "1.Acid-bases chemistry/Buffers: dissociation constant, titration curves, buffering capacity. Ionisation equilibria of acids and bases in aqueous solutions2.Amino acids: structure of the amino acids, properties of amino acids side chains, classes and nomenclature, acid-bases properties, stereochemistry, Modified amino acids.3.Proteins: primary structure, peptides and the peptide bond. Protein purification, solubility, chromatography, electrophoresis. Proteins sequences.4.Proteins: three-dimensional structure. Secondary structure, tertiary structure, quaternary structure. Dynamics of protein structure:  folding and stability.5.Protein Function: oxygen-binding proteins: myoglobin and hemoglobin, immunoglobulins.6.Enzymes: Biological catalysis. Nomenclature and classification of enzymes, specificity of the substrates, co-factors and co-enzymes. Energy of activation and coordinator of reaction. Progress curves. Effect of temperature and pH on the rate of enzymatic activity. 7.Enzymatic kinetics, inhibition and regulation: Effect of the substrate concentration on the rate of enzymatic reactions. Enzyme Inhibition. Reversible inhibition. Models of competitive, non-competitive and simple competitive inhibition. Allosteric regulation.8.Carbohydrates: Monosaccharides, stereoisomerism, classification, configuration and conformation. Derivatives of monosaccharides, oligosaccharides, structural polysaccharides: cellulose and chitin, storage polysaccharides: starch and glycogen, Glycoconjugates: Proteoglycans, Glycoproteins, and Glycolipids.9.Lipids: Classification, fatty acids, triacylglycerols, waxes, glycerophospholipids, esphingolipids, cholesterol.10.Biological membranes. Proteins in membranes: integral and peripheral, the fluid mosaic models, the asymmetry of membranes, transport across membranes: the thermodynamics of transport, passive and active transport.11.Nucleotides and nucleic acids: structure and function of nucleotides. Primary structure of nucleic acids, nucleic acid sequencing, restriction endonucleases. Genome sequency. DNA amplification by the polymerase chain reaction."
$txtAvaliacaoMetodo = "A avaliação será feita por meio de provas escritas."
$txtNotaFinal = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3"
$txtRecuperacao = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$txtLOT2002 = "LOT2002 -  Biologia Celular  (Requisito fraco)`n"
$txtLOT2059 = "LOT2059 -  Química Orgânica Fundamental  (Requisito fraco)`n"

# ---- Row 10: Objetivos content replaced with the professor string ----
$ws.Range("B10").Value = $txtAdriane
$ws.Range("C10").Value = $txtAdriane

# ---- Row 13 (was blank-A / "427823..." B-C): becomes "Programa resumido:" / "Semestral" ----
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---- Row 14: "Short syllabus:" + english short syllabus text (shifted up from row 15) ----
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $txtShortSyllabusPt
$ws.Range("C14").Value = $txtShortSyllabusPt

# ---- Row 15: "Programa:" / "01/01/2018", height grows to 120 ----
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# ---- Row 16: "Syllabus:" / full english syllabus (unchanged content, kept as-is) ----
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $txtSyllabusEn
$ws.Range("C16").Value = $txtSyllabusEn

# ---- Row 17: becomes "Avaliacao:" label only; clear B/C and the custom height ----
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(12).RowHeight

# ---- Row 18: "Metodo:" / professor string (new B/C content), height 60 ----
$ws.Range("A18").Value = "Método:"
$ws.Range("B16").Copy($ws.Range("B18"))
$ws.Range("C16").Copy($ws.Range("C18"))
$ws.Range("B18").Value = $txtAdriane
$ws.Range("C18").Value = $txtAdriane
$ws.Rows.Item(18).RowHeight = 60

# ---- Row 19: "Criterio:" / avaliacao method text ----
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = $txtAvaliacaoMetodo
$ws.Range("C19").Value = $txtAvaliacaoMetodo

# ---- Row 20: "Norma de recuperacao:" / nota final formula ----
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = $txtNotaFinal
$ws.Range("C20").Value = $txtNotaFinal

# ---- Row 21: "Bibliografia:" / recuperacao text, height grows to 120 ----
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = $txtRecuperacao
$ws.Range("C21").Value = $txtRecuperacao
$ws.Rows.Item(21).RowHeight = 120

# ---- Row 22: "Requisitos:" label only; clear B/C and the custom height ----
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Rows.Item(22).RowHeight = $ws.Rows.Item(12).RowHeight

# ---- Row 23: clear old "A23" label, add LOT2002 requirement text to B/C ----
$ws.Range("A23").Value = ""
$ws.Range("B24").Copy($ws.Range("B23"))
$ws.Range("C24").Copy($ws.Range("C23"))
$ws.Range("B23").Value = $txtLOT2002
$ws.Range("C23").Value = $txtLOT2002
$ws.Rows.Item(23).RowHeight = 30

# ---- Row 24: LOT2059 requirement text (shifted up from row 25) ----
$ws.Range("B24").Value = $txtLOT2059
$ws.Range("C24").Value = $txtLOT2059

# ---- Remove the now-redundant last row (content moved up into row 24) ----
$ws.Rows.Item(25).Delete()

